$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1576
$ws.Range("I2").Value = 4288
$ws.Range("J2").Value = 17967
$ws.Range("K2").Value = 94
$ws.Range("L2").Value = 4950
$ws.Range("M2").Value = 285
$ws.Range("N2").Value = 3261
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 73
$ws.Range("Q2").Value = 23
$ws.Range("R2").Value = 243
$ws.Range("S2").Value = 1916
$ws.Range("T2").Value = 3050
$ws.Range("U2").Value = 220
$ws.Range("V2").Value = 28066
$ws.Range("W2").Value = 6
$ws.Range("X2").Value = 27843
$ws.Range("Y2").Value = 50
$ws.Range("Z2").Value = 396
$ws.Range("AA2").Value = 168
